$wb = $excel.ActiveWorkbook
$win = $wb.Windows.Item(1)
try { $win.FirstSheet = 2 ; Write-Output "set ok" } catch { Write-Output "set fail: $_" }
